# Rework the "10A1" math-scores sheet into a "Lop10A1" student-roster sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/sheet-tab.
$ws.Name = "Lop10A1"

# Drop the old score rows (5-11); the new roster only needs rows 1-4.
$ws.Range("A5:D11").Clear()

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Họ và tên"
$ws.Range("B1").Value = "Ngày Sinh"
$ws.Range("C1").Value = "Giới Tính"
$ws.Range("D1").Value = "Địa Chỉ"

# --- Data rows --------------------------------------------------------
# Column B holds birth dates written as literal text (e.g. "1990-04-24"),
# not Excel date serials, so force the cells to Text before writing and
# strip the formatting back off afterwards so no date-parsing happens.
$ws.Range("B2:B4").NumberFormat = "@"

$ws.Range("A2").Value = "Nguyễn Văn An"
$ws.Range("B2").Value = "1990-04-24"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "243 Khuất Duy Tiến-Thanh Xuân-Hà Nội"

$ws.Range("A3").Value = "Nguyễn Thi Thạch Anh"
$ws.Range("B3").Value = "1990-07-18"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "Ba Vì-Hà Tây"

$ws.Range("A4").Value = "Cao Quyết Thắng"
$ws.Range("B4").Value = "1990-03-23"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "Tây Hòa-Phú Yên"

$ws.Range("B2:B4").ClearFormats()
